$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$edits = @(
    @("M28", 211.5),
    @("K28", 273.5),
    @("I28", 273.5),
    @("H28", 316.2353),
    @("I29", 321.5),
    @("L29", 8228.571599999999),
    @("H29", 1625.3077),
    @("K29", 964.5),
    @("N29", -8790.571599999999),
    @("J29", 2742.8572),
    @("M29", -683.5),
    @("J41", 344.8889),
    @("I41", 1433.4445),
    @("M41", -993.4445000000001),
    @("K41", 1433.4445),
    @("N41", -1224.8889),
    @("H41", 889.1667),
    @("L41", 344.8889),
    @("H42", 178.90909),
    @("J42", 260.42856),
    @("K42", 108.75),
    @("M42", 121.25),
    @("N42", -1241.28568),
    @("I42", 36.25),
    @("L42", 781.28568),
    @("N51", -2718.3334),
    @("I51", 800),
    @("L51", 1750.3334),
    @("M51", -316),
    @("H51", 1614.5714),
    @("K51", 800),
    @("J51", 1750.3334),
    @("H61", 165),
    @("K61", 495),
    @("M61", -323),
    @("I61", 165),
    @("J62", 166671180),
    @("I62", 83335020),
    @("M62", -83334396),
    @("N62", -166672428),
    @("L62", 166671180),
    @("H62", 111113740),
    @("K62", 83335020),
    @("N64", -3517.75),
    @("H64", 3029.3784),
    @("L64", 3021.75),
    @("J64", 3021.75),
    @("M64", -2805.111),
    @("I64", 3053.111),
    @("K64", 3053.111),
    @("J65", 166671180),
    @("K65", 416675100),
    @("L65", 833355900),
    @("H65", 111113740),
    @("N65", -833362140),
    @("I65", 83335020),
    @("M65", -416671980),
    @("N67", -4737.75),
    @("J67", 3021.75),
    @("L67", 3021.75),
    @("I67", 3053.111),
    @("M67", -2195.111),
    @("K67", 3053.111),
    @("H67", 3029.3784),
    @("K76", 4653),
    @("H76", 5130.5713),
    @("L76", 5321.6),
    @("M76", -4338),
    @("J76", 5321.6),
    @("I76", 4653),
    @("N76", -5951.6),
    @("I79", 4653),
    @("H79", 5130.5713),
    @("N79", -7505.6),
    @("K79", 4653),
    @("J79", 5321.6),
    @("L79", 5321.6),
    @("M79", -3561),
    @("L92", 1126),
    @("K92", 18518714),
    @("H92", 12820995),
    @("M92", -18517466),
    @("J92", 1126),
    @("N92", -3622),
    @("I92", 18518714),
    @("N98", -3511.8),
    @("H98", 439.5625),
    @("L98", 515.8),
    @("J98", 515.8),
    @("I116", 7659677),
    @("K116", 7659677),
    @("M116", -7656235),
    @("H116", 5904878.5),
    @("L116", 2374.5454),
    @("N116", -9258.545399999999),
    @("J116", 2374.5454),
    @("L122", 1547.4),
    @("H122", 439.5625),
    @("N122", -6447.4),
    @("J122", 515.8),
    @("M132", -9657.9614),
    @("H132", 4097.407),
    @("I132", 4062.6538),
    @("K132", 12187.9614),
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$edits = @(
    @("H2", 1019.13336),
    @("M2", -783.4545000000001),
    @("I2", 896.4545000000001),
    @("K2", 896.4545000000001),
    @("M110", -7504.23),
    @("J110", 3000),
    @("K110", 9549.23),
    @("L110", 3000),
    @("H110", 8676),
    @("N110", -7090),
    @("I110", 9549.23),
    @("I116", 896.4545000000001),
    @("K116", 896.4545000000001),
    @("M116", 1397.5455),
    @("H116", 1019.13336),
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$edits = @(
    @("K3", 896.4545000000001),
    @("M3", -782.4545000000001),
    @("H3", 1019.13336),
    @("I3", 896.4545000000001),
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$edits = @(
    @("J16", 3351.25),
    @("H16", 41669384),
    @("L16", 3351.25),
    @("K16", 125001450),
    @("M16", -125001163),
    @("I16", 125001450),
    @("N16", -3925.25),
    @("N94", -2852),
    @("J94", 1950),
    @("L94", 1950),
    @("H94", 111112810),
    @("M94", -333334079),
    @("I94", 333334530),
    @("K94", 333334530),
    @("K99", 3577039.5),
    @("I99", 3577039.5),
    @("J99", 5671.3335),
    @("N99", -8667.333500000001),
    @("M99", -3575541.5),
    @("H99", 2237776.5),
    @("L99", 5671.3335),
    @("H113", 41669384),
    @("M113", -124999280),
    @("N113", -7691.25),
    @("I113", 125001450),
    @("K113", 125001450),
    @("J113", 3351.25),
    @("L113", 3351.25),
    @("K122", 15386805),
    @("L122", 7581.333),
    @("H122", 3031768.2),
    @("N122", -12481.333),
    @("I122", 5128935),
    @("J122", 2527.111),
    @("M122", -15384355),
    @("K126", 10731118.5),
    @("J126", 5671.3335),
    @("I126", 3577039.5),
    @("N126", -21954.0005),
    @("H126", 2237776.5),
    @("M126", -10728648.5),
    @("L126", 17014.0005),
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$edits = @(
    @("L92", 546.75),
    @("H92", 212),
    @("J92", 182.25),
    @("N92", -3042.75),
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$edits = @(
    @("J80", 2817.375),
    @("L80", 2817.375),
    @("H80", 2847.65),
    @("M80", -1970.75),
    @("N80", -4813.375),
    @("K80", 2968.75),
    @("I80", 2968.75),
    @("M83", -9851.75),
    @("H83", 2847.65),
    @("L83", 14086.875),
    @("I83", 2968.75),
    @("N83", -24070.875),
    @("J83", 2817.375),
    @("K83", 14843.75),
    @("H97", 1147.5),
    @("I97", 1147.5),
    @("M97", -651.5),
    @("K97", 1147.5),
    @("J102", 2700),
    @("M102", -242),
    @("N102", -5944),
    @("H102", 1947.6),
    @("I102", 1864),
    @("L102", 2700),
    @("K102", 1864),
    @("K122", 7917.75),
    @("L122", 5593.9998),
    @("H122", 2229.1765),
    @("N122", -10493.9998),
    @("I122", 2639.25),
    @("J122", 1864.6666),
    @("M122", -5467.75),
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$edits = @(
    @("J16", 4800),
    @("H16", 1254.8182),
    @("L16", 4800),
    @("K16", 900.3),
    @("M16", -730.3),
    @("I16", 900.3),
    @("N16", -5140),
    @("L40", 2166.6667),
    @("H40", 1787.75),
    @("M40", -1424.4),
    @("K40", 1560.4),
    @("N40", -2438.6667),
    @("J40", 2166.6667),
    @("I40", 1560.4),
    @("M93", 65.07690000000002),
    @("K93", 1182.9231),
    @("N93", -4521),
    @("I93", 1182.9231),
    @("J93", 2025),
    @("L93", 2025),
    @("H93", 1381.0588),
    @("K122", 5100),
    @("L122", 8634.706200000001),
    @("H122", 2610.4546),
    @("N122", -13534.7062),
    @("I122", 1700),
    @("J122", 2878.2354),
    @("M122", -2650),
    @("L136", 10499.4),
    @("N136", -15599.4),
    @("J136", 3499.8),
    @("H136", 5441894.5),
    @("K136", 20857678.5),
    @("M136", -20855128.5),
    @("I136", 6952559.5),
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$edits = @(
    @("L54", 16974.75),
    @("H54", 16974.75),
    @("J54", 16974.75),
    @("N54", -18014.75),
    @("M107", 18),
    @("I107", 634),
    @("J107", 1222),
    @("L107", 3666),
    @("N107", -7506),
    @("K107", 1902),
    @("H107", 794.36365),
    @("H113", 20834428),
    @("M113", -107144114),
    @("I113", 35715428),
    @("K113", 107146284),
    @("K122", 214289430),
    @("L122", 6784.6158),
    @("H122", 37038770),
    @("N122", -11684.6158),
    @("I122", 71429810),
    @("J122", 2261.5386),
    @("M122", -214286980),
    @("M132", -5227.625),
    @("H132", 2522.8667),
    @("J132", 2450.8572),
    @("L132", 7352.571599999999),
    @("N132", -12412.5716),
    @("I132", 2585.875),
    @("K132", 7757.625),
)
foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}
